{"js": "// Split the single run of text in the Title, Author and Abstract\n// paragraphs into one run per word, with the separating spaces kept\n// as their own runs (matching the target OOXML diff).\n\nfunction wrapBodyOoxml(innerBodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    innerBodyXml +\n    \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\nfunction escapeXmlText(s) {\n  return s.replace(/&/g, \"&amp;\").replace(/</g, \"&lt;\").replace(/>/g, \"&gt;\");\n}\n\n// Rebuild `paragraph` so that its text is split word-by-word, each word\n// (and each separating space) becoming its own <w:r>, while preserving\n// the paragraph's own formatting (<w:pPr>) and the formatting of its\n// (single, unformatted-in-practice) run (<w:rPr>), taken from the\n// paragraph's current OOXML.\nasync function splitIntoWordRuns(paragraph) {\n  paragraph.load(\"text\");\n  const ooxml = paragraph.getOoxml();\n  await context.sync();\n\n  const text = paragraph.text;\n  // Nothing to split.\n  if (!text || text.indexOf(\" \") === -1) {\n    return;\n  }\n\n  const fullXml = ooxml.value;\n  const bodyIdx = fullXml.indexOf(\"<w:body>\");\n  const afterBody = bodyIdx >= 0 ? fullXml.slice(bodyIdx) : fullXml;\n  const pMatch = afterBody.match(/<w:p[ >][\\s\\S]*?<\\/w:p>/);\n  const pXml = pMatch ? pMatch[0] : \"<w:p/>\";\n\n  const pPrMatch = pXml.match(/<w:pPr>[\\s\\S]*?<\\/w:pPr>/);\n  const pPrXml = pPrMatch ? pPrMatch[0] : \"\";\n\n  const rPrMatch = pXml.match(/<w:rPr>[\\s\\S]*?<\\/w:rPr>/);\n  const rPrXml = rPrMatch ? rPrMatch[0] : \"\";\n\n  // Keep the separating spaces as their own tokens.\n  const tokens = text.split(/( )/).filter((tok) => tok.length > 0);\n\n  const runsXml = tokens\n    .map(\n      (tok) =>\n        \"<w:r>\" +\n        rPrXml +\n        '<w:t xml:space=\"preserve\">' +\n        escapeXmlText(tok) +\n        \"</w:t></w:r>\"\n    )\n    .join(\"\");\n\n  const newParagraphXml = \"<w:p>\" + pPrXml + runsXml + \"</w:p>\";\n  paragraph.insertOoxml(wrapBodyOoxml(newParagraphXml), \"Replace\");\n  await context.sync();\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the Title, Author and Abstract paragraphs by their current text\n// (robust to exact index) rather than assuming fixed positions.\nlet titlePara = null;\nlet authorPara = null;\nlet abstractPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text === \"Factsheet: Laws of indices\") {\n    titlePara = p;\n  } else if (p.text === \"Tom Coleman\") {\n    authorPara = p;\n  } else if (p.text === \"A list of laws of indices.\") {\n    abstractPara = p;\n  }\n}\n\nif (titlePara) {\n  await splitIntoWordRuns(titlePara);\n}\nif (authorPara) {\n  await splitIntoWordRuns(authorPara);\n}\nif (abstractPara) {\n  await splitIntoWordRuns(abstractPara);\n}\n", "ps1": "# Split the single run of text in the Title, Author and Abstract\n# paragraphs into one run per word, with the separating spaces kept\n# as their own runs (matching the target OOXML diff).\n\nfunction New-WrappedOoxml($innerBodyXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $innerBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n}\n\nfunction ConvertTo-XmlEscaped($s) {\n    return $s.Replace('&', '&amp;').Replace('<', '&lt;').Replace('>', '&gt;')\n}\n\n# Rebuild $paragraph so that its text is split word-by-word, each word\n# (and each separating space) becoming its own run, while preserving the\n# paragraph's own formatting (<w:pPr>) and the formatting of its\n# (single, unformatted-in-practice) run (<w:rPr>), taken from the\n# paragraph's current OOXML.\nfunction Split-ParagraphIntoWordRuns($paragraph) {\n    $r = $paragraph.Range\n    $text = $r.Text\n    # Paragraph Range.Text includes the trailing paragraph mark; strip it.\n    $text = $text.Substring(0, $text.Length - 1)\n\n    if ($text -notmatch ' ') {\n        return\n    }\n\n    $xml = $r.WordOpenXML\n    $bodyIdx = $xml.IndexOf(\"<w:body>\")\n    $after = $xml.Substring($bodyIdx)\n\n    $pXml = ''\n    if ($after -match '(?s)<w:p[ >].*?</w:p>') {\n        $pXml = $matches[0]\n    }\n\n    $pPrXml = ''\n    if ($pXml -match '(?s)<w:pPr>.*?</w:pPr>') {\n        $pPrXml = $matches[0]\n    }\n\n    $rPrXml = ''\n    if ($pXml -match '(?s)<w:rPr>.*?</w:rPr>') {\n        $rPrXml = $matches[0]\n    }\n\n    $words = $text -split ' '\n    $tokens = @()\n    for ($i = 0; $i -lt $words.Length; $i++) {\n        if ($i -gt 0) {\n            $tokens += ' '\n        }\n        $tokens += $words[$i]\n    }\n\n    $runsXml = ''\n    foreach ($tok in $tokens) {\n        $runsXml += '<w:r>' + $rPrXml + '<w:t xml:space=\"preserve\">' + (ConvertTo-XmlEscaped $tok) + '</w:t></w:r>'\n    }\n\n    $newParaXml = '<w:p>' + $pPrXml + $runsXml + '</w:p>'\n    [void]$r.InsertXML((New-WrappedOoxml $newParaXml))\n}\n\n$d = $word.ActiveDocument\n\n# Locate the Title, Author and Abstract paragraphs by their current text\n# (robust to exact index) rather than assuming fixed positions.\n$targets = @()\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    $t = $t.Substring(0, $t.Length - 1)\n    if ($t -eq \"Factsheet: Laws of indices\" -or $t -eq \"Tom Coleman\" -or $t -eq \"A list of laws of indices.\") {\n        $targets += $p\n    }\n}\n\nforeach ($p in $targets) {\n    Split-ParagraphIntoWordRuns $p\n}\n"}
